# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect freshly scraped counts, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Row => New value for column F
$updates = @{
    3  = 3278
    4  = 239
    5  = 134
    7  = 1709
    8  = 1644
    10 = 378
    14 = 33
    16 = 240
    18 = 9
    24 = 236
    25 = 109
    29 = 347
    30 = 2239
    31 = 12
    33 = 476
    34 = 447
    35 = 568
    40 = 531
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
